# Update countries & provincias Spain
# Applies the 27-May-2020 refresh: updated timestamp, updated case counts for
# several countries, and corrects the row order for a couple of country pairs
# whose names/data had been mismatched (Banglades/Bielorrusia/Ecuador and
# Groenlandia/Islas Turcas y Caicos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp refresh (row 1) ---
$ws.Cells.Item(1,1).Value = "Datos actualizados a 27 de Mayo de 2020 a las 11:05"

# --- Estados Unidos (row 4) ---
$ws.Cells.Item(4,2).Value = 1725278
$ws.Cells.Item(4,3).Value = 3
$ws.Cells.Item(4,4).Value = 479973
$ws.Cells.Item(4,5).Value = 1144733

# --- Alemania (row 11) ---
$ws.Cells.Item(11,4).Value = 162800
$ws.Cells.Item(11,5).Value = 9990

# --- Belgica (row 22) ---
$ws.Cells.Item(22,2).Value = 57592
$ws.Cells.Item(22,3).Value = 137
$ws.Cells.Item(22,4).Value = 15465
$ws.Cells.Item(22,5).Value = 32763
$ws.Cells.Item(22,7).Value = 30
$ws.Cells.Item(22,8).Value = 9364

# --- Row 25 becomes Banglades (fresh data), row 26 becomes Bielorrusia
#     (prior row-25 data), row 27 becomes Ecuador (prior row-26 data) ---
$ws.Cells.Item(26,1).Value = "Bielorrusia"
$ws.Cells.Item(26,2).Value = 38059
$ws.Cells.Item(26,4).Value = 15086
$ws.Cells.Item(26,5).Value = 22765
$ws.Cells.Item(26,8).Value = 208

$ws.Cells.Item(27,1).Value = "Ecuador"
$ws.Cells.Item(27,2).Value = 37355
$ws.Cells.Item(27,4).Value = 18003
$ws.Cells.Item(27,5).Value = 16149
$ws.Cells.Item(27,8).Value = 3203

$ws.Cells.Item(25,1).Value = "Banglades"
$ws.Cells.Item(25,2).Value = 38292
$ws.Cells.Item(25,3).Value = 1541
$ws.Cells.Item(25,4).Value = 7925
$ws.Cells.Item(25,5).Value = 29823
$ws.Cells.Item(25,7).Value = 22
$ws.Cells.Item(25,8).Value = 544

# --- Indonesia (row 35) ---
$ws.Cells.Item(35,2).Value = 23851
$ws.Cells.Item(35,3).Value = 686
$ws.Cells.Item(35,4).Value = 6057
$ws.Cells.Item(35,5).Value = 16321
$ws.Cells.Item(35,7).Value = 55
$ws.Cells.Item(35,8).Value = 1473

# --- Filipinas (row 46) ---
$ws.Cells.Item(46,2).Value = 15049
$ws.Cells.Item(46,3).Value = 380
$ws.Cells.Item(46,4).Value = 3506
$ws.Cells.Item(46,5).Value = 10639
$ws.Cells.Item(46,7).Value = 18
$ws.Cells.Item(46,8).Value = 904

# --- Barein (row 53) ---
$ws.Cells.Item(53,5).Value = 4413
$ws.Cells.Item(53,7).Value = 1
$ws.Cells.Item(53,8).Value = 15

# --- Finlandia (row 67) ---
$ws.Cells.Item(67,2).Value = 6692
$ws.Cells.Item(67,3).Value = 64
$ws.Cells.Item(67,5).Value = 1280

# --- Sri Lanka (row 103) ---
$ws.Cells.Item(103,4).Value = 732
$ws.Cells.Item(103,5).Value = 577

# --- Malta (row 133) ---
$ws.Cells.Item(133,5).Value = 119
$ws.Cells.Item(133,7).Value = 1
$ws.Cells.Item(133,8).Value = 7

# --- Row 207 becomes Groenlandia, row 208 becomes Islas Turcas y Caicos
#     (the two rows' labels/data swap) ---
$ws.Cells.Item(207,1).Value = "Groenlandia"
$ws.Cells.Item(207,4).Value = 11
$ws.Cells.Item(207,8).Value = 0

$ws.Cells.Item(208,1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(208,4).Value = 10
$ws.Cells.Item(208,8).Value = 1
